$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "studentid"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "secrettoken"

# Row 2 - Kassim Balogun
$ws.Range("A2").Value = 34010231
$ws.Range("B2").Value = "Kassim Balogun"
$ws.Range("C2").Value = "er2345"

# Row 3 - Richard Alipui
$ws.Range("A3").Value = 34562212
$ws.Range("B3").Value = "Richard Alipui"
$ws.Range("C3").Value = "tf4534"

# Row 4 - Papa Kofi Gyekye
$ws.Range("A4").Value = 34516066
$ws.Range("B4").Value = "Papa Kofi Gyekye"
$ws.Range("C4").Value = "2hg432"

# Row 5 - Peter Perez
$ws.Range("A5").Value = 34917894
$ws.Range("B5").Value = "Peter Perez"
$ws.Range("C5").Value = "6ty435"

# Column widths widened to fit the new longer name / token columns
$ws.Columns("B").ColumnWidth = 13.3
$ws.Columns("C").ColumnWidth = 9.8

# View state: zoomed in, selection moved to B4
$excel.ActiveWindow.Zoom = 156
$ws.Range("B4").Select() | Out-Null

Write-Output "done"
